# Automatische test-sync: 2025-06-24 20:08:50
#
# Adds the new "Korting voor wederverkopers?" mail-log entry (row 17) to the
# Logs sheet, extends the conditional-formatting ranges to cover it, and
# re-sorts the Dashboard category/count table to reflect the updated tally
# (Offerte / Prijsaanvraag now has 2 entries).

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new log row ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "Korting voor wederverkopers?"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Biedt u speciale prijzen voor wederverkopers?"
$logs.Range("D17").Value = "Offerte / Prijsaanvraag"
$logs.Range("F17").Value = "2025-06-24 20:08:18"
$logs.Range("G17").Value = "Nee"

# Widen the conditional-formatting sqref (D2:D16 -> D2:D17, G2:G16 -> G2:G17)
# so the new row is covered by the category / "Beantwoord" colour rules.
$logs.Range("D2:D16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D17"))
$logs.Range("G2:G16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G17"))

# --- Dashboard sheet: re-sort category counts ---
# "Offerte / Prijsaanvraag" moves up to row 5 (now tied at 2), the other two
# categories shift down a row, and the Offerte count increments to 2.
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("A6").Value = "Productinformatie"
$dash.Range("A7").Value = "IT / Technisch probleem"
$dash.Range("B7").Value = 2
